$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy date-column formatting (style index 2) to the two newly appended rows (114, 115)
$ws.Range("A113").Copy()
$ws.Range("A114:A115").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update data rows 90-115 with the refreshed report values
$ws.Cells.Item(90, 1).Value = 44232
$ws.Cells.Item(90, 2).Value = 28
$ws.Cells.Item(90, 3).Value = 122
$ws.Cells.Item(90, 4).Value = 168.5804695380619

$ws.Cells.Item(91, 1).Value = 44233
$ws.Cells.Item(91, 2).Value = 28
$ws.Cells.Item(91, 3).Value = 123
$ws.Cells.Item(91, 4).Value = 169.962276665423

$ws.Cells.Item(92, 1).Value = 44234
$ws.Cells.Item(92, 2).Value = 19
$ws.Cells.Item(92, 3).Value = 127
$ws.Cells.Item(92, 4).Value = 175.4895051748677

$ws.Cells.Item(93, 1).Value = 44235
$ws.Cells.Item(93, 2).Value = 14
$ws.Cells.Item(93, 3).Value = 137
$ws.Cells.Item(93, 4).Value = 189.3075764484793

$ws.Cells.Item(94, 1).Value = 44236
$ws.Cells.Item(94, 2).Value = 28
$ws.Cells.Item(94, 3).Value = 119
$ws.Cells.Item(94, 4).Value = 164.4350481559784

$ws.Cells.Item(95, 1).Value = 44237
$ws.Cells.Item(95, 2).Value = 7
$ws.Cells.Item(95, 3).Value = 102
$ws.Cells.Item(95, 4).Value = 140.9443269908386

$ws.Cells.Item(96, 1).Value = 44238
$ws.Cells.Item(96, 2).Value = 13
$ws.Cells.Item(96, 3).Value = 103
$ws.Cells.Item(96, 4).Value = 142.3261341181998

$ws.Cells.Item(97, 1).Value = 44239
$ws.Cells.Item(97, 2).Value = 10
$ws.Cells.Item(97, 3).Value = 108
$ws.Cells.Item(97, 4).Value = 149.2351697550056

$ws.Cells.Item(98, 1).Value = 44240
$ws.Cells.Item(98, 2).Value = 11
$ws.Cells.Item(98, 3).Value = 93
$ws.Cells.Item(98, 4).Value = 128.5080628445882

$ws.Cells.Item(99, 1).Value = 44241
$ws.Cells.Item(99, 2).Value = 20
$ws.Cells.Item(99, 3).Value = 88
$ws.Cells.Item(99, 4).Value = 121.5990272077823

$ws.Cells.Item(100, 1).Value = 44242
$ws.Cells.Item(100, 2).Value = 19
$ws.Cells.Item(100, 3).Value = 100
$ws.Cells.Item(100, 4).Value = 138.1807127361163

$ws.Cells.Item(101, 1).Value = 44243
$ws.Cells.Item(101, 2).Value = 13
$ws.Cells.Item(101, 3).Value = 117
$ws.Cells.Item(101, 4).Value = 161.6714339012561

$ws.Cells.Item(102, 1).Value = 44244
$ws.Cells.Item(102, 2).Value = 2
$ws.Cells.Item(102, 3).Value = 126
$ws.Cells.Item(102, 4).Value = 174.1076980475065

$ws.Cells.Item(103, 1).Value = 44245
$ws.Cells.Item(103, 2).Value = 25
$ws.Cells.Item(103, 3).Value = 132
$ws.Cells.Item(103, 4).Value = 182.3985408116735

$ws.Cells.Item(104, 1).Value = 44246
$ws.Cells.Item(104, 2).Value = 27
$ws.Cells.Item(104, 3).Value = 136
$ws.Cells.Item(104, 4).Value = 187.9257693211182

$ws.Cells.Item(105, 1).Value = 44247
$ws.Cells.Item(105, 2).Value = 20
$ws.Cells.Item(105, 3).Value = 154
$ws.Cells.Item(105, 4).Value = 212.7982976136191

$ws.Cells.Item(106, 1).Value = 44248
$ws.Cells.Item(106, 2).Value = 26
$ws.Cells.Item(106, 3).Value = 158
$ws.Cells.Item(106, 4).Value = 218.3255261230637

$ws.Cells.Item(107, 1).Value = 44249
$ws.Cells.Item(107, 2).Value = 23
$ws.Cells.Item(107, 3).Value = 156
$ws.Cells.Item(107, 4).Value = 215.5619118683414

$ws.Cells.Item(108, 1).Value = 44250
$ws.Cells.Item(108, 2).Value = 31
$ws.Cells.Item(108, 3).Value = 180
$ws.Cells.Item(108, 4).Value = 248.7252829250093

$ws.Cells.Item(109, 1).Value = 44251
$ws.Cells.Item(109, 2).Value = 6
$ws.Cells.Item(109, 3).Value = 181
$ws.Cells.Item(109, 4).Value = 250.1070900523705

$ws.Cells.Item(110, 1).Value = 44252
$ws.Cells.Item(110, 2).Value = 23
$ws.Cells.Item(110, 3).Value = 184
$ws.Cells.Item(110, 4).Value = 254.252511434454

$ws.Cells.Item(111, 1).Value = 44253
$ws.Cells.Item(111, 2).Value = 51
$ws.Cells.Item(111, 3).Value = 206
$ws.Cells.Item(111, 4).Value = 284.6522682363996

$ws.Cells.Item(112, 1).Value = 44254
$ws.Cells.Item(112, 2).Value = 21
$ws.Cells.Item(112, 3).Value = 191
$ws.Cells.Item(112, 4).Value = 263.9251613259821

$ws.Cells.Item(113, 1).Value = 44255
$ws.Cells.Item(113, 2).Value = 29

$ws.Cells.Item(114, 1).Value = 44256
$ws.Cells.Item(114, 2).Value = 45

$ws.Cells.Item(115, 1).Value = 44257
$ws.Cells.Item(115, 2).Value = 16

